# CodeSystem-cgm-system-cs.xlsx: FHIR IG terminology/metadata correction
#
# - Metadata sheet: "Experimental" row (row 7) value B7 was blank; set it to
#   the literal text "false".
# - Metadata sheet: "Date" row (row 8) value B8 is bumped to the new
#   generation timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Write B7 with a leading apostrophe so Excel stores it as TEXT ("false")
# rather than coercing it to the boolean FALSE.
$ws.Range("B7").Value = "'false"

# Re-apply B6's cell format to B7 via PasteSpecial(Formats) so the text entry
# above keeps the sheet's normal (non "quote-prefixed") cell style, matching
# every other data row.
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Update the generation date/time string.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
